$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet colour-codes assignment rows by status (green = submitted-ish,
# red = further out, gold = approaching deadline). Re-use existing formats
# via copy/paste-special so the underlying theme-coloured fills are shared
# instead of spawning duplicate raw-RGB fills.

# Row 15 (Data Visualization, Part 1): gold -> green
$ws.Range("A15:D15").Interior.Color = 5287936

# Row 16 (Data Visualization, Part 2 -- "ass2 done"): red -> gold (copy format from row 14,
# since the gold fill is theme-based and can't be reproduced bit-for-bit from a raw RGB value)
$ws.Range("A14:D14").Copy()
$ws.Range("A16:D16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 17 (Embedded Image Processing, assignment 4): gold -> green
$ws.Range("A17:D17").Interior.Color = 5287936

# Row 19 C cell: re-apply the same red fill so the style is canonicalised
# (drops a redundant applyFont flag it had picked up previously)
$ws.Range("C19").Interior.Color = 255

# New row 20: Embedded Image Processing, assignment 5, due 14-Mar-2023, not submitted
$ws.Range("A20").Value = "Embedded Image Processing "
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = 44999
$ws.Range("D20").Value = "No"

# New light-gold fill (theme accent colour, lighter tint than the existing gold rows)
$ws.Range("A20:D20").Interior.Color = 13431551
$ws.Range("C20").NumberFormat = "d-mmm"

# Update view: scrolled down a bit, with F18 selected
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("F18").Select()
